$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.138.68'
$ws.Range('E2').Value = '  +5.71%  '
$ws.Range('D3').Value = '2.418.38'
$ws.Range('E3').Value = '  +2.02%  '
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.86%  '
$ws.Range('D5').Value = '573.94'
$ws.Range('E5').Value = '  +2.72%  '
$ws.Range('D6').Value = '146.39'
$ws.Range('E6').Value = '  +6.47%  '
$ws.Range('E8').Value = '  +2.17%  '
$ws.Range('D9').Value = '2.458.80'
$ws.Range('E9').Value = '  +3.92%  '
$ws.Range('E10').Value = '  +6.05%  '
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('D12').Value = '5.24'
$ws.Range('E12').Value = '  +3.10%  '
$ws.Range('E13').Value = '  +4.82%  '
$ws.Range('D14').Value = '27.42'
$ws.Range('E14').Value = '  +7.21%  '
$ws.Range('D15').Value = '0.0000178'
$ws.Range('E15').Value = '  +8.31%  '
$ws.Range('D16').Value = '2.858.22'
$ws.Range('D17').Value = '62.992.45'
$ws.Range('E17').Value = '  +5.48%  '
$ws.Range('D18').Value = '2.444.03'
$ws.Range('E18').Value = '  +3.09%  '
$ws.Range('D19').Value = '7.92'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('D20').Value = '11.01'
$ws.Range('E20').Value = '  +5.05%  '
$ws.Range('D21').Value = '328.86'
$ws.Range('E21').Value = '  +2.47%  '
$ws.Range('E22').Value = '  +2.34%  '
$ws.Range('D23').Value = '2.07'
$ws.Range('E23').Value = '  +14.62%  '
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.20%  '
$ws.Range('D25').Value = '65.70'
$ws.Range('E25').Value = '  +2.55%  '
$ws.Range('D26').Value = '633.42'
$ws.Range('E26').Value = '  +13.66%  '
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').Value = '8.56'
$ws.Range('E27').Value = '  +5.24%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '0.0₃0988'
$ws.Range('E28').Value = '  +7.70%  '
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '2.537.87'
$ws.Range('E29').Value = '  +2.13%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '8.23'
$ws.Range('E30').Value = '  +3.11%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = '1.42'
$ws.Range('E31').Value = '  +9.18%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = '0.138'
$ws.Range('E32').Value = '  +6.22%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.86'
$ws.Range('E33').Value = '  +4.64%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '1.50'
$ws.Range('E34').Value = '  +4.99%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').Value = '0.995'
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = '4.76'
$ws.Range('E36').Value = '  +5.15%  '
$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').Value = '0.374'
$ws.Range('E37').Value = '  +2.29%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = '153.06'
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').Value = '5.43'
$ws.Range('E39').Value = '  +9.15%  '
$ws.Range('B40').Value = 'EthereumClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D40').Value = '18.70'
$ws.Range('E40').Value = '  +3.12%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = '2.75'
$ws.Range('E41').Value = '  +14.22%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '1.78'
$ws.Range('E42').Value = '  +8.74%  '
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D44').Value = '0.0₆0296'
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '145.52'
$ws.Range('E45').Value = '  +4.89%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '3.60'
$ws.Range('E46').Value = '  +2.69%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '20.56'
$ws.Range('E47').Value = '  +8.31%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = '0.603'
$ws.Range('E48').Value = '  +3.21%  '
$ws.Range('B49').Value = 'Hedera'
$ws.Range('C49').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D49').Value = '0.0517'
$ws.Range('E49').Value = '  +3.82%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').Value = '12.68'
$ws.Range('E50').Value = '  +8.58%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '0.0918'
$ws.Range('E51').Value = '  +2.19%  '
